$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row at position 31, shifting existing rows 31-71 down to 32-72
$ws.Rows("31:31").Insert()

# Materialize all cells A31:Y31 (matching the fully-populated empty cells of
# surrounding rows) without introducing a new cell style.
$ws.Range("A31:Y31").Style = "Normal"

# Populate the new row 31 with the new entry (September_Details / September_Date)
$ws.Cells.Item(31, 18).Value = "transfer"
$ws.Cells.Item(31, 19).Value = "2024-09-05 16:25:07"
